$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Bgn"
$ws.Cells.Item(2,3).Value = "Tlr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 38.718679
$ws.Cells.Item(2,8).Value = 116.156037
$ws.Cells.Item(2,9).Value = 0.01404461724059496
$ws.Cells.Item(2,10).Value = 0.01404461724059496
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 44.04866033333334
$ws.Cells.Item(2,14).Value = 132.145981
$ws.Cells.Item(2,15).Value = 0.3636142564479216
$ws.Cells.Item(2,16).Value = 0.3636142564479216
$ws.Cells.Item(2,17).Value = 1705.505939826367
$ws.Cells.Item(2,18).Value = 15349.5534584373
$ws.Cells.Item(2,19).Value = 0.005106823055034598
$ws.Cells.Item(2,20).Value = 0.005106823055034598

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Bgn"
$ws.Cells.Item(3,3).Value = "Tlr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 38.718679
$ws.Cells.Item(3,8).Value = 116.156037
$ws.Cells.Item(3,9).Value = 0.01404461724059496
$ws.Cells.Item(3,10).Value = 0.01404461724059496
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.226320666666667
$ws.Cells.Item(3,14).Value = 6.678962
$ws.Cells.Item(3,15).Value = 0.01837790134135009
$ws.Cells.Item(3,16).Value = 0.01837790134135009
$ws.Cells.Item(3,17).Value = 86.20019524373267
$ws.Cells.Item(3,18).Value = 775.801757193594
$ws.Cells.Item(3,19).Value = 0.0002581105900246787
$ws.Cells.Item(3,20).Value = 0.0002581105900246787

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Bgn"
$ws.Cells.Item(4,3).Value = "Tlr2"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 38.718679
$ws.Cells.Item(4,8).Value = 116.156037
$ws.Cells.Item(4,9).Value = 0.01404461724059496
$ws.Cells.Item(4,10).Value = 0.01404461724059496
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 74.86619966666667
$ws.Cells.Item(4,14).Value = 224.598599
$ws.Cells.Item(4,15).Value = 0.6180078422107282
$ws.Cells.Item(4,16).Value = 0.6180078422107284
$ws.Cells.Item(4,17).Value = 2898.720352843574
$ws.Cells.Item(4,18).Value = 26088.48317559217
$ws.Cells.Item(4,19).Value = 0.008679683595535686
$ws.Cells.Item(4,20).Value = 0.008679683595535687

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Bgn"
$ws.Cells.Item(5,3).Value = "Tlr2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2511.398112333333
$ws.Cells.Item(5,8).Value = 7534.194336999999
$ws.Cells.Item(5,9).Value = 0.9109718135392577
$ws.Cells.Item(5,10).Value = 0.9109718135392579
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 44.04866033333334
$ws.Cells.Item(5,14).Value = 132.145981
$ws.Cells.Item(5,15).Value = 0.3636142564479216
$ws.Cells.Item(5,16).Value = 0.3636142564479216
$ws.Cells.Item(5,17).Value = 110623.7224119455
$ws.Cells.Item(5,18).Value = 995613.5017075095
$ws.Cells.Item(5,19).Value = 0.3312423386250918
$ws.Cells.Item(5,20).Value = 0.3312423386250919

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Bgn"
$ws.Cells.Item(6,3).Value = "Tlr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2511.398112333333
$ws.Cells.Item(6,8).Value = 7534.194336999999
$ws.Cells.Item(6,9).Value = 0.9109718135392577
$ws.Cells.Item(6,10).Value = 0.9109718135392579
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.226320666666667
$ws.Cells.Item(6,14).Value = 6.678962
$ws.Cells.Item(6,15).Value = 0.01837790134135009
$ws.Cells.Item(6,16).Value = 0.01837790134135009
$ws.Cells.Item(6,17).Value = 5591.177519715354
$ws.Cells.Item(6,18).Value = 50320.59767743819
$ws.Cells.Item(6,19).Value = 0.01674175011397525
$ws.Cells.Item(6,20).Value = 0.01674175011397525

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Bgn"
$ws.Cells.Item(7,3).Value = "Tlr2"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2511.398112333333
$ws.Cells.Item(7,8).Value = 7534.194336999999
$ws.Cells.Item(7,9).Value = 0.9109718135392577
$ws.Cells.Item(7,10).Value = 0.9109718135392579
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 74.86619966666667
$ws.Cells.Item(7,14).Value = 224.598599
$ws.Cells.Item(7,15).Value = 0.6180078422107282
$ws.Cells.Item(7,16).Value = 0.6180078422107284
$ws.Cells.Item(7,17).Value = 188018.8325204371
$ws.Cells.Item(7,18).Value = 1692169.492683934
$ws.Cells.Item(7,19).Value = 0.5629877248001905
$ws.Cells.Item(7,20).Value = 0.5629877248001908

$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Bgn"
$ws.Cells.Item(8,3).Value = "Tlr2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 5.092134333333333
$ws.Cells.Item(8,8).Value = 15.276403
$ws.Cells.Item(8,9).Value = 0.001847094980935658
$ws.Cells.Item(8,10).Value = 0.001847094980935659
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 44.04866033333334
$ws.Cells.Item(8,14).Value = 132.145981
$ws.Cells.Item(8,15).Value = 0.3636142564479216
$ws.Cells.Item(8,16).Value = 0.3636142564479216
$ws.Cells.Item(8,17).Value = 224.3016956207048
$ws.Cells.Item(8,18).Value = 2018.715260586343
$ws.Cells.Item(8,19).Value = 0.0006716300680816073
$ws.Cells.Item(8,20).Value = 0.0006716300680816073

$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Bgn"
$ws.Cells.Item(9,3).Value = "Tlr2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 5.092134333333333
$ws.Cells.Item(9,8).Value = 15.276403
$ws.Cells.Item(9,9).Value = 0.001847094980935658
$ws.Cells.Item(9,10).Value = 0.001847094980935659
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.226320666666667
$ws.Cells.Item(9,14).Value = 6.678962
$ws.Cells.Item(9,15).Value = 0.01837790134135009
$ws.Cells.Item(9,16).Value = 0.01837790134135009
$ws.Cells.Item(9,17).Value = 11.33672390374289
$ws.Cells.Item(9,18).Value = 102.030515133686
$ws.Cells.Item(9,19).Value = 0.00003394572932773845
$ws.Cells.Item(9,20).Value = 0.00003394572932773845

$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Bgn"
$ws.Cells.Item(10,3).Value = "Tlr2"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 5.092134333333333
$ws.Cells.Item(10,8).Value = 15.276403
$ws.Cells.Item(10,9).Value = 0.001847094980935658
$ws.Cells.Item(10,10).Value = 0.001847094980935659
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 74.86619966666667
$ws.Cells.Item(10,14).Value = 224.598599
$ws.Cells.Item(10,15).Value = 0.6180078422107282
$ws.Cells.Item(10,16).Value = 0.6180078422107284
$ws.Cells.Item(10,17).Value = 381.2287457288219
$ws.Cells.Item(10,18).Value = 3431.058711559398
$ws.Cells.Item(10,19).Value = 0.001141519183526312
$ws.Cells.Item(10,20).Value = 0.001141519183526313

$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Bgn"
$ws.Cells.Item(11,3).Value = "Tlr2"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 201.6251223333333
$ws.Cells.Item(11,8).Value = 604.875367
$ws.Cells.Item(11,9).Value = 0.07313647423921157
$ws.Cells.Item(11,10).Value = 0.07313647423921157
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 44.04866033333334
$ws.Cells.Item(11,14).Value = 132.145981
$ws.Cells.Item(11,15).Value = 0.3636142564479216
$ws.Cells.Item(11,16).Value = 0.3636142564479216
$ws.Cells.Item(11,17).Value = 8881.316528327781
$ws.Cells.Item(11,18).Value = 79931.84875495003
$ws.Cells.Item(11,19).Value = 0.02659346469971349
$ws.Cells.Item(11,20).Value = 0.02659346469971349

$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Bgn"
$ws.Cells.Item(12,3).Value = "Tlr2"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 201.6251223333333
$ws.Cells.Item(12,8).Value = 604.875367
$ws.Cells.Item(12,9).Value = 0.07313647423921157
$ws.Cells.Item(12,10).Value = 0.07313647423921157
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.226320666666667
$ws.Cells.Item(12,14).Value = 6.678962
$ws.Cells.Item(12,15).Value = 0.01837790134135009
$ws.Cells.Item(12,16).Value = 0.01837790134135009
$ws.Cells.Item(12,17).Value = 448.8821767698949
$ws.Cells.Item(12,18).Value = 4039.939590929054
$ws.Cells.Item(12,19).Value = 0.001344094908022422
$ws.Cells.Item(12,20).Value = 0.001344094908022422

$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Bgn"
$ws.Cells.Item(13,3).Value = "Tlr2"
$ws.Cells.Item(13,4).Value = "M2"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 201.6251223333333
$ws.Cells.Item(13,8).Value = 604.875367
$ws.Cells.Item(13,9).Value = 0.07313647423921157
$ws.Cells.Item(13,10).Value = 0.07313647423921157
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 74.86619966666667
$ws.Cells.Item(13,14).Value = 224.598599
$ws.Cells.Item(13,15).Value = 0.6180078422107282
$ws.Cells.Item(13,16).Value = 0.6180078422107284
$ws.Cells.Item(13,17).Value = 15094.90666642343
$ws.Cells.Item(13,18).Value = 135854.1599978108
$ws.Cells.Item(13,19).Value = 0.04519891463147566
$ws.Cells.Item(13,20).Value = 0.04519891463147566
